$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.282.50'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '1.708.91'
$ws.Range('E3').Value = '  -0.66%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '223.78'
$ws.Range('E5').Value = '  -2.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5286'
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2633'
$ws.Range('E8').Value = '  -4.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06533'
$ws.Range('E9').Value = '  -3.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.85'
$ws.Range('E10').Value = '  -2.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07634'
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.566'
$ws.Range('E12').Value = '  -3.22%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.947.69'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.684.80'
$ws.Range('E14').Value = '  -1.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5735'
$ws.Range('E15').Value = '  -3.63%  '
$ws.Range('D16').Value = '0.0₅8169'
$ws.Range('E16').Value = '  -2.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.12'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('D18').Value = '27.276.35'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.71'
$ws.Range('E19').Value = '  +3.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  +0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.680'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.46'
$ws.Range('E22').Value = '  -3.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.943'
$ws.Range('E23').Value = '  -4.41%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.005'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.32'
$ws.Range('E25').Value = '  -2.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.751'
$ws.Range('E26').Value = '  +7.14%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1215'
$ws.Range('E27').Value = '  -2.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.249'
$ws.Range('E28').Value = '  -2.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '16.27'
$ws.Range('E29').Value = '  -3.58%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05362'
$ws.Range('E30').Value = '  -4.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.487'
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.406'
$ws.Range('E33').Value = '  -2.82%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.636'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.870'
$ws.Range('E35').Value = '  +1.14%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.425'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9494'
$ws.Range('E37').Value = '  -2.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5862'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01627'
$ws.Range('E39').Value = '  -0.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.867'
$ws.Range('E40').Value = '  -0.13%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8391'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.036.73'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.98'
$ws.Range('E44').Value = '  -1.15%  '
$ws.Range('D45').Value = '1.853.51'
$ws.Range('E45').Value = '  -0.38%  '
$ws.Range('E46').Value = '  +5.18%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.92'
$ws.Range('E47').Value = '  -2.85%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4492'
$ws.Range('E48').Value = '  +2.61%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.002'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.084'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05242'
$ws.Range('E51').Value = '  -0.45%  '
